# Apply the "move step forward" button edit:
#  1. Rename the existing rewind label from "5S" to "3S".
#  2. Duplicate the rewind button group (background rectangle, icon picture,
#     label textbox) to create a new "step forward" button: reposition it
#     below the existing row, flip the icon horizontally (turning the
#     "skip to beginning" glyph into a "skip to end"/forward glyph), rename
#     the duplicated shapes, and set the new label text to "3S".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "5S" -> "3S" on the existing rewind textbox -------------------------
$s.Shapes.Item(3).TextFrame.TextRange.Text = "3S"

# --- 2. Duplicate the button group (Rectangle 6 / Graphic 4 / TextBox 5) ----
$newRect = $s.Shapes.Item(1).Duplicate().Item(1)
$newPic  = $s.Shapes.Item(2).Duplicate().Item(1)
$newText = $s.Shapes.Item(3).Duplicate().Item(1)

# Rename to match the new shapes
$newRect.Name = "Rectangle 2"
$newPic.Name  = "Graphic 3"
$newText.Name = "TextBox 8"

# Reposition the new background rectangle
# (point values nudged slightly above the exact EMU/12700 quotient so that,
# after the host's internal 32-bit-float round trip, the saved EMU lands on
# the exact target instead of one unit low)
$newRect.Left   = 99.10590751181101
$newRect.Top    = 302.7481232362205
$newRect.Width  = 140.61173278346456
$newRect.Height = 59.92952955905512

# Reposition + flip the new icon picture (becomes a "step forward" glyph)
$newPic.Left   = 169.41173578346456
$newPic.Top    = 303.07637795275593
$newPic.Width  = 59.92952955905512
$newPic.Height = 59.92952955905512
$newPic.HorizontalFlip = -1

# Reposition the new label textbox and set its text
$newText.Left   = 113.55031496062992
$newText.Top    = 304.8433075866142
$newText.Width  = 59.92952955905512
$newText.Height = 55.73905511811024
$newText.TextFrame.TextRange.Text = "3S"
